$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Insert 3 new rows in the middle of the existing data block (rows 4:6),
# so the existing SUM formulas (originally SUM(B2:B4) etc.) expand
# automatically to cover the new rows, turning into SUM(B2:B7) etc.
$ws.Rows("4:6").Insert()

# Re-write the six data rows (3 pre-existing days shifted down, plus
# 3 new days inserted above them) in descending date order.
$data = @(
  @(2, 45449, 1,   1,  0, 0),
  @(3, 45448, 141, 27, 3, 111),
  @(4, 45447, 0,   0,  0, 0),
  @(5, 45446, 212, 1,  0, 211),
  @(6, 45445, 252, 26, 9, 217),
  @(7, 45444, 209, 30, 6, 173)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]
  $ws.Cells.Item($r, 5).Value = $row[5]
}

# Update the active selection to match the new totals row.
$ws.Range("B8:E8").Select()
